$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.073001804238757018,
    -0.0099999996804811531,
    -0.0089999996841818586,
    0.28399386409065741,
    -0.0059999996942483591,
    -0.0059999996827073687,
    -0.019999999621591158,
    0.015865885145244896,
    -0.0059999996736923578,
    -0.0059999996685533574,
    -0.0044999996743229076,
    -0.0059999996664918953,
    -0.0059999996609834128,
    -0.01199999963357623,
    -0.0059999996581412418,
    -0.0059999996569790603,
    -0.005999999655483812,
    -0.0089999996422251982,
    -0.0089999996844243313,
    -0.024772725637410176,
    -0.0089999996794372095,
    -0.008999999679126347,
    -0.0089999996814045247,
    -0.041999999534531263,
    -0.04199999953206035,
    -0.0059999996813928647,
    -0.0059999996796884503,
    -0.015363491073277835,
    -0.011999999641648884,
    -0.019999999604948915,
    -0.014999999623217519,
    -0.020999999596855723,
    -0.0059999996611930229
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
